# Update column F (dSF) values per re-pulled / re-calculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 4
    4  = 1
    5  = 2
    6  = 5
    8  = 1
    9  = -3
    11 = 4
    12 = -2
    13 = -2
    14 = -3
    15 = 2
    16 = 4
    17 = 2
    18 = 3
    19 = 2
    20 = 4
    21 = -3
    22 = -3
    24 = -1
    25 = 2
    26 = 5
    28 = -3
    29 = -1
    30 = 3
    31 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
